$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 테이블목록
$ws2 = $wb.Worksheets.Item(2)   # 컬럼정의서

# --- Fill in the new rows on the "테이블목록" sheet ---
# SEQ column
$ws1.Cells.Item(2,1).Value = 1
$ws1.Cells.Item(3,1).Value = 2
$ws1.Cells.Item(4,1).Value = 3

# 테이블ID column
$ws1.Cells.Item(2,2).Value = "dept"
$ws1.Cells.Item(3,2).Value = "employee"
$ws1.Cells.Item(4,2).Value = "company"

# 테이블명 column
$ws1.Cells.Item(2,3).Value = "부서"
$ws1.Cells.Item(3,3).Value = "직원"
$ws1.Cells.Item(4,3).Value = "회사"

# 엔터티정의 column (values start with a literal apostrophe, so double the
# leading apostrophe: the first is consumed as Excel's quote-prefix marker,
# the second is kept as a literal character in the cell text)
$ws1.Cells.Item(2,4).Value = "''부서' 속성을 갖는 엔터티"
$ws1.Cells.Item(3,4).Value = "''직원' 속성을 갖는 엔터티"
$ws1.Cells.Item(4,4).Value = "''회사' 속성을 갖는 엔터티"

# --- Widen column D on "테이블목록" ---
$ws1.Columns.Item(4).ColumnWidth = 21.571428571428573

# --- Switch the active sheet from 컬럼정의서 to 테이블목록 ---
$ws2.Activate()
$ws2.Range("G19").Select() | Out-Null

$ws1.Activate()
$ws1.Range("C12").Select() | Out-Null
